$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32 (pushes old rows 32..51 down to 33..52)
$ws.Rows("32:32").Insert()

# Fill the brand-new row 32 with its data
$ws.Cells.Item(32, 1).Value = 10
$ws.Cells.Item(32, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(32, 3).Value = "La Araucanía"
$ws.Cells.Item(32, 4).Value = 44567
$ws.Cells.Item(32, 5).Value = 9
$ws.Cells.Item(32, 6).Value = 100112030
$ws.Cells.Item(32, 7).Value = "Poroto granado"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 130
$ws.Cells.Item(32, 11).Value = 27000
$ws.Cells.Item(32, 12).Value = 28000
$ws.Cells.Item(32, 13).Value = 27615
$ws.Cells.Item(32, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(32, 15).Value = "Región del Maule"
$ws.Cells.Item(32, 16).Value = 1105
$ws.Cells.Item(32, 17).Value = 25
$ws.Cells.Item(32, 18).Value = "Hortaliza"

# Insert a new row at position 52 (pushes old row 51, now sitting at 52, down to 53)
$ws.Rows("52:52").Insert()

# Fill the brand-new row 52 with its data
$ws.Cells.Item(52, 1).Value = 10
$ws.Cells.Item(52, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(52, 3).Value = "La Araucanía"
$ws.Cells.Item(52, 4).Value = 44568
$ws.Cells.Item(52, 5).Value = 9
$ws.Cells.Item(52, 6).Value = 100112030
$ws.Cells.Item(52, 7).Value = "Poroto granado"
$ws.Cells.Item(52, 8).Value = "Sin especificar"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 80
$ws.Cells.Item(52, 11).Value = 30000
$ws.Cells.Item(52, 12).Value = 30000
$ws.Cells.Item(52, 13).Value = 30000
$ws.Cells.Item(52, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(52, 15).Value = "Región del Maule"
$ws.Cells.Item(52, 16).Value = 1200
$ws.Cells.Item(52, 17).Value = 25
$ws.Cells.Item(52, 18).Value = "Hortaliza"
